# Updates the "Recorded By" (column G) values for specific rows to match
# the re-ordered author list produced by the upstream sync job.
# Only the rows listed below changed; other rows with identical text were
# left untouched in the source diff, so we must target rows explicitly
# rather than doing a blanket text transform.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$gUpdates = @{
    "2" = 'backup@backdoor.com, system, System'
    "3" = 'dnasr281@gmail.com, System'
    "4" = 'backup@backdoor.com, System'
    "5" = 'backup@backdoor.com, System'
    "6" = 'dnasr281@gmail.com, System'
    "8" = 'backup@backdoor.com, System'
    "10" = 'dnasr281@gmail.com, System'
    "12" = 'dnasr281@gmail.com, System'
    "13" = 'dnasr281@gmail.com, System'
    "14" = 'dnasr281@gmail.com, System'
    "15" = 'dnasr281@gmail.com, System'
    "18" = 'dnasr281@gmail.com, System'
    "19" = 'dnasr281@gmail.com, System'
    "20" = 'dnasr281@gmail.com, System'
    "21" = 'dnasr281@gmail.com, System'
    "22" = 'dnasr281@gmail.com, System'
    "24" = 'dnasr281@gmail.com, System'
    "29" = 'backup@backdoor.com, system, System'
    "30" = 'dnasr281@gmail.com, System'
    "31" = 'backup@backdoor.com, System'
    "32" = 'backup@backdoor.com, System'
    "33" = 'dnasr281@gmail.com, System'
    "35" = 'backup@backdoor.com, System'
    "37" = 'dnasr281@gmail.com, System'
    "39" = 'dnasr281@gmail.com, System'
    "40" = 'dnasr281@gmail.com, System'
    "41" = 'dnasr281@gmail.com, System'
    "42" = 'dnasr281@gmail.com, System'
    "45" = 'dnasr281@gmail.com, System'
    "46" = 'dnasr281@gmail.com, System'
    "47" = 'dnasr281@gmail.com, System'
    "48" = 'dnasr281@gmail.com, System'
    "49" = 'dnasr281@gmail.com, System'
    "51" = 'dnasr281@gmail.com, System'
    "56" = 'backup@backdoor.com, system, System'
    "57" = 'dnasr281@gmail.com, System'
    "58" = 'backup@backdoor.com, System'
    "59" = 'backup@backdoor.com, System'
    "60" = 'dnasr281@gmail.com, System'
    "62" = 'backup@backdoor.com, System'
    "64" = 'dnasr281@gmail.com, System'
    "66" = 'dnasr281@gmail.com, System'
    "67" = 'dnasr281@gmail.com, System'
    "68" = 'dnasr281@gmail.com, System'
    "69" = 'dnasr281@gmail.com, System'
    "72" = 'dnasr281@gmail.com, System'
    "73" = 'dnasr281@gmail.com, System'
    "74" = 'dnasr281@gmail.com, System'
    "75" = 'dnasr281@gmail.com, System'
    "76" = 'dnasr281@gmail.com, System'
    "78" = 'dnasr281@gmail.com, System'
    "83" = 'backup@backdoor.com, System'
    "84" = 'backup@backdoor.com, System'
    "85" = 'backup@backdoor.com, System'
    "86" = 'dnasr281@gmail.com, System'
    "87" = 'dnasr281@gmail.com, System'
    "88" = 'dnasr281@gmail.com, System'
    "89" = 'dnasr281@gmail.com, System'
    "93" = 'dnasr281@gmail.com, System'
    "95" = 'dnasr281@gmail.com, System'
    "102" = 'dnasr281@gmail.com, System'
    "109" = 'backup@backdoor.com, System'
    "110" = 'backup@backdoor.com, System'
    "111" = 'backup@backdoor.com, System'
    "112" = 'dnasr281@gmail.com, System'
    "113" = 'dnasr281@gmail.com, System'
    "114" = 'dnasr281@gmail.com, System'
    "115" = 'dnasr281@gmail.com, System'
    "119" = 'dnasr281@gmail.com, System'
    "121" = 'dnasr281@gmail.com, System'
    "128" = 'dnasr281@gmail.com, System'
    "135" = 'backup@backdoor.com, System'
    "136" = 'backup@backdoor.com, System'
    "137" = 'backup@backdoor.com, System'
    "138" = 'dnasr281@gmail.com, System'
    "139" = 'dnasr281@gmail.com, System'
    "140" = 'dnasr281@gmail.com, System'
    "141" = 'dnasr281@gmail.com, System'
    "145" = 'dnasr281@gmail.com, System'
    "147" = 'dnasr281@gmail.com, System'
    "154" = 'dnasr281@gmail.com, System'
}

foreach ($row in $gUpdates.Keys) {
    $ws.Range("G" + $row).Value = $gUpdates[$row]
}
